$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("C9").Value = -0.0224
$ws.Range("D9").Value = 0.0448
$ws.Range("E9").Value = 0.0063
$ws.Range("F9").Value = -0.003
$ws.Range("G9").Value = -0.0062
$ws.Range("H9").Value = -0.1304
$ws.Range("I9").Value = -0.0761
$ws.Range("J9").Value = -0.0439
$ws.Range("K9").Value = -0.0496
$ws.Range("L9").Value = -0.0213
$ws.Range("M9").Value = -0.0314

# Row 10
$ws.Range("C10").Value = 0.1231
$ws.Range("D10").Value = 0.1382
$ws.Range("E10").Value = 0.2464
$ws.Range("F10").Value = 0.2121
$ws.Range("G10").Value = 0.1649
$ws.Range("H10").Value = 0.1224
$ws.Range("I10").Value = -0.1438
$ws.Range("J10").Value = -0.1364
$ws.Range("K10").Value = -0.1463
$ws.Range("L10").Value = -0.1453

# Row 12
$ws.Range("C12").Value = -0.1684
$ws.Range("D12").Value = 0.0324
$ws.Range("E12").Value = 0.638
$ws.Range("F12").Value = -0.7926
$ws.Range("G12").Value = 0.0168
$ws.Range("H12").Value = -0.7905

# Row 15
$ws.Range("C15").Value = -2.2358
$ws.Range("D15").Value = -2.3933
$ws.Range("E15").Value = -1.7935
$ws.Range("F15").Value = -1.9824
$ws.Range("G15").Value = -3.3881
$ws.Range("H15").Value = -2.838
$ws.Range("I15").Value = -2.1902
$ws.Range("J15").Value = -1.8154
$ws.Range("K15").Value = -2.96
$ws.Range("L15").Value = -1.2147
$ws.Range("M15").Value = -0.628

# Row 36
$ws.Range("C36").Value = 0.0002
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0.0002
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0.0002
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0

# Row 37
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0

# Row 39
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0.165
$ws.Range("F39").Value = 0.1584
$ws.Range("G39").Value = 0.1549
$ws.Range("H39").Value = 0.1512

# Row 42
$ws.Range("C42").Value = 0.0001
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0.1949
$ws.Range("F42").Value = -0.0343
$ws.Range("G42").Value = -0.0773
$ws.Range("H42").Value = -0.0785
$ws.Range("I42").Value = 0.0957
$ws.Range("J42").Value = 0.0641
$ws.Range("K42").Value = 0.0376
$ws.Range("L42").Value = 0.0193
$ws.Range("M42").Value = -0.1131
